# The commit removes the stray "TEst -1" sub-title line (and the line
# break that introduced it) from the title placeholder on slide 1, so the
# title box reads:
#   Deep Learning Specialisation
#   TOPICS COVERED
#   (blank line)
# instead of:
#   Deep Learning Specialisation
#   TOPICS COVERED
#   TEst -1

$p  = $ppt.ActivePresentation
$s  = $p.Slides.Item(1)
$sh = $s.Shapes.Item(1)
$tr = $sh.TextFrame.TextRange

$full = $tr.Text
$idx0 = $full.IndexOf("TEst -1")

if ($idx0 -ge 0) {
    # $idx0 (0-based) lands exactly on the line-break character that sits
    # right before the "TEst -1" run when addressed through the 1-based
    # Characters() indexer, so this range covers "<br>TEst -1" (8 chars)
    # and removes both the run and the now-superfluous break in one go.
    $victim = $tr.Characters($idx0, 8)
    $victim.Text = ""
}
